$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need forced Text format
# so Excel does not normalize them (stripping trailing zeros, etc.)

$ws.Range("D2").Value = "60.488.18"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "2.637.64"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.08"
$ws.Range("E5").Value = "  +5.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.53"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +5.28%  "
$ws.Range("D9").Value = "2.663.63"
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.84"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("E12").Value = "  +7.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").Value = "3.105.91"
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "60.438.97"
$ws.Range("E15").Value = "  +3.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.97"
$ws.Range("E16").Value = "  +5.71%  "
$ws.Range("E17").Value = "  +4.54%  "
$ws.Range("D18").Value = "2.659.61"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  +3.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.96"
$ws.Range("E20").Value = "  +2.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +2.96%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.09"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.442"
$ws.Range("E25").Value = "  +5.04%  "
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("E28").Value = "  +4.38%  "
$ws.Range("D29").Value = "0.0₃0806"
$ws.Range("E29").Value = "  +9.77%  "
$ws.Range("E31").Value = "  +4.26%  "
$ws.Range("E32").Value = "  +4.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.21"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.903"
$ws.Range("E36").Value = "  +8.68%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.16"
$ws.Range("E37").Value = "  +5.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.899"
$ws.Range("E38").Value = "  +9.43%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  +6.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "302.74"
$ws.Range("E41").Value = "  +6.82%  "
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0990"
$ws.Range("E43").Value = "  +4.87%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.994"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.604"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "128.50"
$ws.Range("E47").Value = "  +14.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.31"
$ws.Range("E48").Value = "  +1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.70"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +5.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.64"
$ws.Range("E51").Value = "  +3.94%  "
